# Update the monitoring data: re-sort the "empadronador" table by
# total_registros (descending) and refresh the totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ALBIRENA GARCIA ANGEELO ALONSO", 139),
    @("MANUEL LEUNARDO PRADO BAILON", 124),
    @("ALAMA NIMA CLARITZA MABEL", 105),
    @("URRIOLA ARISMENDIZ INGRID MARYURI", 104),
    @("CORDOVA CARMEN ANGIE NATALLY", 102),
    @("AGURTO ORDINOLA LISBET JAQUELIN", 96),
    @("VEGA ROBLEDO FERNANDO ERNESTO", 88),
    @("ROMAN GALECIO MARITZA DEL PILAR", 85),
    @("MARYURI OJEDA VALLE", 80),
    @("JUAREZ CARMEN PIERRE ALEXANDER", 78),
    @("ATOCHE PALACIOS LUIS ANGEL", 73),
    @("CARREÑO PALACIOS KATHERINE DE LOS MILAGROS", 69),
    @("RUIDIAS FRIAS MELISSA VICTORIA", 64),
    @("BERNAOLA CARMEN ZUMIKO YASHURY", 61)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
